$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" shared string to "In Translation" everywhere it's used.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"

# Narrow columns E/F on Overview sheet and C on zh-cn/de-de sheets.
$overview.Range("E1").ColumnWidth = 13.4101845877511
$overview.Range("F1").ColumnWidth = 13.4101845877511

$zhcn.Range("C1").ColumnWidth = 13.4101845877511
$dede.Range("C1").ColumnWidth = 13.4101845877511
